$wb = $excel.ActiveWorkbook

# ===== Sheet: In Report =====
$ws = $wb.Worksheets.Item("In Report")
$ws.Rows.Item(5).Insert()
$ws.Range("A4").Copy()
$ws.Range("A5:X5").PasteSpecial(-4122)
$ws.Range("A5").Value = "SL #"
$ws.Range("B5").Value = "Id"
$ws.Range("C5").Value = "Container Number"
$ws.Range("D5").Value = "Container Size"
$ws.Range("E5").Value = "Container Type"
$ws.Range("F5").Value = "Current Depot Unit"
$ws.Range("G5").Value = "Permitted Depot Unit"
$ws.Range("H5").Value = "Agent"
$ws.Range("I5").Value = "Mlo"
$ws.Range("J5").Value = "Source Location"
$ws.Range("K5").Value = "Import Vessel Name"
$ws.Range("L5").Value = "Import Rotation Number"
$ws.Range("M5").Value = "Gate In Date"
$ws.Range("N5").Value = "Container Condition Name"
$ws.Range("O5").Value = "Di Agent"
$ws.Range("P5").Value = "Di Mlo"
$ws.Range("Q5").Value = "Di Date"
$ws.Range("R5").Value = "Remarks"
$ws.Range("S5").Value = "Damage Area Name"
$ws.Range("T5").Value = "Damage Part Name"
$ws.Range("U5").Value = "Damage Description"
$ws.Range("V5").Value = "Damage Component"
$ws.Range("W5").Value = "Damage Type"
$ws.Range("X5").Value = "Repair Type"
$ws.Columns.Item(2).ColumnWidth = 6.142857142857143
$ws.Columns.Item(2).Hidden = $true
$ws.Columns.Item(3).ColumnWidth = 24.142857142857142
$ws.Columns.Item(4).ColumnWidth = 18.714285714285715
$ws.Columns.Item(5).ColumnWidth = 22.285714285714285
$ws.Columns.Item(6).ColumnWidth = 26.0
$ws.Columns.Item(7).ColumnWidth = 27.714285714285715
$ws.Columns.Item(8).ColumnWidth = 11.571428571428571
$ws.Columns.Item(9).ColumnWidth = 8.0
$ws.Columns.Item(10).ColumnWidth = 18.714285714285715
$ws.Columns.Item(11).ColumnWidth = 24.142857142857142
$ws.Columns.Item(12).ColumnWidth = 29.571428571428573
$ws.Columns.Item(13).ColumnWidth = 17.0
$ws.Columns.Item(14).ColumnWidth = 33.142857142857146
$ws.Columns.Item(15).ColumnWidth = 15.142857142857142
$ws.Columns.Item(16).ColumnWidth = 11.571428571428571
$ws.Columns.Item(17).ColumnWidth = 11.571428571428571
$ws.Columns.Item(18).ColumnWidth = 9.714285714285714
$ws.Columns.Item(19).ColumnWidth = 24.142857142857142
$ws.Columns.Item(20).ColumnWidth = 22.285714285714285
$ws.Columns.Item(21).ColumnWidth = 22.285714285714285
$ws.Columns.Item(22).ColumnWidth = 27.714285714285715
$ws.Columns.Item(23).ColumnWidth = 20.571428571428573
$ws.Columns.Item(24).ColumnWidth = 18.714285714285715

# ===== Sheet: Out Empty Report =====
$ws = $wb.Worksheets.Item("Out Empty Report")
$ws.Rows.Item(5).Insert()
$ws.Range("A4").Copy()
$ws.Range("A5:AE5").PasteSpecial(-4122)
$ws.Range("A5").Value = "SL #"
$ws.Range("B5").Value = "Id"
$ws.Range("C5").Value = "Container Number"
$ws.Range("D5").Value = "Container Size"
$ws.Range("E5").Value = "Container Type"
$ws.Range("F5").Value = "Current Depot Unit"
$ws.Range("G5").Value = "Permitted Depot Unit"
$ws.Range("H5").Value = "Agent"
$ws.Range("I5").Value = "Mlo"
$ws.Range("J5").Value = "Source Location"
$ws.Range("K5").Value = "Import Vessel Name"
$ws.Range("L5").Value = "Import Rotation Number"
$ws.Range("M5").Value = "Gate In Date"
$ws.Range("N5").Value = "Container Condition Name"
$ws.Range("O5").Value = "Destination Location"
$ws.Range("P5").Value = "Export Vessel Name"
$ws.Range("Q5").Value = "Export Rotation Number"
$ws.Range("R5").Value = "Stuffing Date"
$ws.Range("S5").Value = "Gate Out Date"
$ws.Range("T5").Value = "Eir Number"
$ws.Range("U5").Value = "Seal No"
$ws.Range("V5").Value = "Commodity"
$ws.Range("W5").Value = "Vat"
$ws.Range("X5").Value = "Cbm"
$ws.Range("Y5").Value = "Weight"
$ws.Range("Z5").Value = "Account"
$ws.Range("AA5").Value = "Forwarder"
$ws.Range("AB5").Value = "Di Agent"
$ws.Range("AC5").Value = "Di Mlo"
$ws.Range("AD5").Value = "Di Date"
$ws.Range("AE5").Value = "Remarks"
$ws.Columns.Item(2).ColumnWidth = 6.142857142857143
$ws.Columns.Item(2).Hidden = $true
$ws.Columns.Item(3).ColumnWidth = 24.142857142857142
$ws.Columns.Item(4).ColumnWidth = 18.714285714285715
$ws.Columns.Item(5).ColumnWidth = 22.285714285714285
$ws.Columns.Item(6).ColumnWidth = 26.0
$ws.Columns.Item(7).ColumnWidth = 27.714285714285715
$ws.Columns.Item(8).ColumnWidth = 11.571428571428571
$ws.Columns.Item(9).ColumnWidth = 8.0
$ws.Columns.Item(10).ColumnWidth = 18.714285714285715
$ws.Columns.Item(11).ColumnWidth = 24.142857142857142
$ws.Columns.Item(12).ColumnWidth = 29.571428571428573
$ws.Columns.Item(13).ColumnWidth = 17.0
$ws.Columns.Item(14).ColumnWidth = 33.142857142857146
$ws.Columns.Item(15).ColumnWidth = 20.571428571428573
$ws.Columns.Item(16).ColumnWidth = 24.142857142857142
$ws.Columns.Item(17).ColumnWidth = 29.571428571428573
$ws.Columns.Item(18).ColumnWidth = 17.0
$ws.Columns.Item(19).ColumnWidth = 18.714285714285715
$ws.Columns.Item(20).ColumnWidth = 17.0
$ws.Columns.Item(21).ColumnWidth = 13.285714285714286
$ws.Columns.Item(22).ColumnWidth = 17.0
$ws.Columns.Item(23).ColumnWidth = 6.142857142857143
$ws.Columns.Item(24).ColumnWidth = 9.714285714285714
$ws.Columns.Item(25).ColumnWidth = 11.571428571428571
$ws.Columns.Item(26).ColumnWidth = 11.571428571428571
$ws.Columns.Item(27).ColumnWidth = 11.571428571428571
$ws.Columns.Item(28).ColumnWidth = 15.142857142857142
$ws.Columns.Item(29).ColumnWidth = 11.571428571428571
$ws.Columns.Item(30).ColumnWidth = 11.571428571428571
$ws.Columns.Item(31).ColumnWidth = 9.714285714285714

# ===== Sheet: Out Laden Report =====
$ws = $wb.Worksheets.Item("Out Laden Report")
$ws.Rows.Item(5).Insert()
$ws.Range("A4").Copy()
$ws.Range("A5:AE5").PasteSpecial(-4122)
$ws.Range("A5").Value = "SL #"
$ws.Range("B5").Value = "Id"
$ws.Range("C5").Value = "Container Number"
$ws.Range("D5").Value = "Container Size"
$ws.Range("E5").Value = "Container Type"
$ws.Range("F5").Value = "Current Depot Unit"
$ws.Range("G5").Value = "Permitted Depot Unit"
$ws.Range("H5").Value = "Agent"
$ws.Range("I5").Value = "Mlo"
$ws.Range("J5").Value = "Source Location"
$ws.Range("K5").Value = "Import Vessel Name"
$ws.Range("L5").Value = "Import Rotation Number"
$ws.Range("M5").Value = "Gate In Date"
$ws.Range("N5").Value = "Container Condition Name"
$ws.Range("O5").Value = "Destination Location"
$ws.Range("P5").Value = "Export Vessel Name"
$ws.Range("Q5").Value = "Export Rotation Number"
$ws.Range("R5").Value = "Stuffing Date"
$ws.Range("S5").Value = "Gate Out Date"
$ws.Range("T5").Value = "Eir Number"
$ws.Range("U5").Value = "Seal No"
$ws.Range("V5").Value = "Commodity"
$ws.Range("W5").Value = "Vat"
$ws.Range("X5").Value = "Cbm"
$ws.Range("Y5").Value = "Weight"
$ws.Range("Z5").Value = "Account"
$ws.Range("AA5").Value = "Forwarder"
$ws.Range("AB5").Value = "Di Agent"
$ws.Range("AC5").Value = "Di Mlo"
$ws.Range("AD5").Value = "Di Date"
$ws.Range("AE5").Value = "Remarks"
$ws.Columns.Item(2).ColumnWidth = 6.142857142857143
$ws.Columns.Item(2).Hidden = $true
$ws.Columns.Item(3).ColumnWidth = 24.142857142857142
$ws.Columns.Item(4).ColumnWidth = 18.714285714285715
$ws.Columns.Item(5).ColumnWidth = 22.285714285714285
$ws.Columns.Item(6).ColumnWidth = 26.0
$ws.Columns.Item(7).ColumnWidth = 27.714285714285715
$ws.Columns.Item(8).ColumnWidth = 11.571428571428571
$ws.Columns.Item(9).ColumnWidth = 8.0
$ws.Columns.Item(10).ColumnWidth = 18.714285714285715
$ws.Columns.Item(11).ColumnWidth = 24.142857142857142
$ws.Columns.Item(12).ColumnWidth = 29.571428571428573
$ws.Columns.Item(13).ColumnWidth = 17.0
$ws.Columns.Item(14).ColumnWidth = 33.142857142857146
$ws.Columns.Item(15).ColumnWidth = 20.571428571428573
$ws.Columns.Item(16).ColumnWidth = 24.142857142857142
$ws.Columns.Item(17).ColumnWidth = 29.571428571428573
$ws.Columns.Item(18).ColumnWidth = 17.0
$ws.Columns.Item(19).ColumnWidth = 18.714285714285715
$ws.Columns.Item(20).ColumnWidth = 17.0
$ws.Columns.Item(21).ColumnWidth = 13.285714285714286
$ws.Columns.Item(22).ColumnWidth = 17.0
$ws.Columns.Item(23).ColumnWidth = 6.142857142857143
$ws.Columns.Item(24).ColumnWidth = 9.714285714285714
$ws.Columns.Item(25).ColumnWidth = 11.571428571428571
$ws.Columns.Item(26).ColumnWidth = 11.571428571428571
$ws.Columns.Item(27).ColumnWidth = 11.571428571428571
$ws.Columns.Item(28).ColumnWidth = 15.142857142857142
$ws.Columns.Item(29).ColumnWidth = 11.571428571428571
$ws.Columns.Item(30).ColumnWidth = 11.571428571428571
$ws.Columns.Item(31).ColumnWidth = 9.714285714285714

# ===== Sheet: Stock Report =====
$ws = $wb.Worksheets.Item("Stock Report")
$ws.Rows.Item(5).Insert()
$ws.Range("A4").Copy()
$ws.Range("A5:AB5").PasteSpecial(-4122)
$ws.Range("A5").Value = "SL #"
$ws.Range("B5").Value = "Id"
$ws.Range("C5").Value = "Container Number"
$ws.Range("D5").Value = "Container Size"
$ws.Range("E5").Value = "Container Type"
$ws.Range("F5").Value = "Agent"
$ws.Range("G5").Value = "Mlo"
$ws.Range("H5").Value = "Current Depot Unit"
$ws.Range("I5").Value = "Permitted Depot Unit"
$ws.Range("J5").Value = "Import Vessel Name"
$ws.Range("K5").Value = "Import Rotation Number"
$ws.Range("L5").Value = "Source Location"
$ws.Range("M5").Value = "Gate In Date"
$ws.Range("N5").Value = "Container Condition Name"
$ws.Range("O5").Value = "Bay Location"
$ws.Range("P5").Value = "Storage Day"
$ws.Range("Q5").Value = "Container Status"
$ws.Range("R5").Value = "Container Status Name"
$ws.Range("S5").Value = "Di Agent"
$ws.Range("T5").Value = "Di Mlo"
$ws.Range("U5").Value = "Di Date"
$ws.Range("V5").Value = "Remarks"
$ws.Range("W5").Value = "Damage Area Name"
$ws.Range("X5").Value = "Damage Part Name"
$ws.Range("Y5").Value = "Damage Description"
$ws.Range("Z5").Value = "Damage Component"
$ws.Range("AA5").Value = "Damage Type"
$ws.Range("AB5").Value = "Repair Type"
$ws.Columns.Item(2).ColumnWidth = 6.142857142857143
$ws.Columns.Item(2).Hidden = $true
$ws.Columns.Item(3).ColumnWidth = 24.142857142857142
$ws.Columns.Item(4).ColumnWidth = 18.714285714285715
$ws.Columns.Item(5).ColumnWidth = 22.285714285714285
$ws.Columns.Item(6).ColumnWidth = 11.571428571428571
$ws.Columns.Item(7).ColumnWidth = 8.0
$ws.Columns.Item(8).ColumnWidth = 26.0
$ws.Columns.Item(9).ColumnWidth = 27.714285714285715
$ws.Columns.Item(10).ColumnWidth = 24.142857142857142
$ws.Columns.Item(11).ColumnWidth = 29.571428571428573
$ws.Columns.Item(12).ColumnWidth = 18.714285714285715
$ws.Columns.Item(13).ColumnWidth = 17.0
$ws.Columns.Item(14).ColumnWidth = 33.142857142857146
$ws.Columns.Item(15).ColumnWidth = 15.142857142857142
$ws.Columns.Item(16).ColumnWidth = 17.0
$ws.Columns.Item(17).ColumnWidth = 18.714285714285715
$ws.Columns.Item(18).ColumnWidth = 26.0
$ws.Columns.Item(19).ColumnWidth = 15.142857142857142
$ws.Columns.Item(20).ColumnWidth = 11.571428571428571
$ws.Columns.Item(21).ColumnWidth = 11.571428571428571
$ws.Columns.Item(22).ColumnWidth = 9.714285714285714
$ws.Columns.Item(23).ColumnWidth = 24.142857142857142
$ws.Columns.Item(24).ColumnWidth = 22.285714285714285
$ws.Columns.Item(25).ColumnWidth = 22.285714285714285
$ws.Columns.Item(26).ColumnWidth = 27.714285714285715
$ws.Columns.Item(27).ColumnWidth = 20.571428571428573
$ws.Columns.Item(28).ColumnWidth = 18.714285714285715
